$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1615.3334
$ws.Range("J17").Value = 1615.3334
$ws.Range("L17").Value = 4846.0002
$ws.Range("N17").Value = -5182.0002

$ws.Range("H19").Value = 1408.9286
$ws.Range("J19").Value = 1560.3334
$ws.Range("L19").Value = 1560.3334
$ws.Range("N19").Value = -1910.3334

$ws.Range("H28").Value = 14749.5
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 14749.5
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 14749.5
$ws.Range("N28").Value = -15719.5
$ws.Range("M28").ClearContents()

$ws.Range("H33").Value = 93.333336
$ws.Range("I33").Value = 93.333336
$ws.Range("K33").Value = 93.333336
$ws.Range("M33").Value = 135.666664

$ws.Range("H62").Value = 5119.75
$ws.Range("I62").Value = 1500
$ws.Range("J62").Value = 6326.3335
$ws.Range("K62").Value = 1500
$ws.Range("L62").Value = 6326.3335
$ws.Range("M62").Value = -876
$ws.Range("N62").Value = -7574.3335

$ws.Range("H65").Value = 5119.75
$ws.Range("I65").Value = 1500
$ws.Range("J65").Value = 6326.3335
$ws.Range("K65").Value = 7500
$ws.Range("L65").Value = 31631.6675
$ws.Range("M65").Value = -4380
$ws.Range("N65").Value = -37871.6675

$ws.Range("H115").Value = 485.6
$ws.Range("I115").Value = 485.6
$ws.Range("K115").Value = 1456.8
$ws.Range("M115").Value = 110.1999999999998

$ws.Range("H121").Value = 997.5
$ws.Range("J121").Value = 997.5
$ws.Range("L121").Value = 2992.5
$ws.Range("N121").Value = -6486.5

$ws.Range("H137").Value = 2064.074
$ws.Range("I137").Value = 1110.6
$ws.Range("K137").Value = 3331.8
$ws.Range("M137").Value = -781.7999999999997

$ws.Range("H138").Value = 2185.658
$ws.Range("J138").Value = 2347.2666
$ws.Range("L138").Value = 7041.7998
$ws.Range("N138").Value = -17321.7998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2227.7317
$ws.Range("I32").Value = 1596.6571
$ws.Range("J32").Value = 5909
$ws.Range("K32").Value = 1596.6571
$ws.Range("L32").Value = 5909
$ws.Range("M32").Value = -1309.6571
$ws.Range("N32").Value = -6483

$ws.Range("H44").Value = 29929.666
$ws.Range("J44").Value = 29929.666
$ws.Range("L44").Value = 29929.666
$ws.Range("N44").Value = -30905.666

$ws.Range("H45").Value = 3002674.2
$ws.Range("J45").Value = 2788.8333
$ws.Range("L45").Value = 2788.8333
$ws.Range("N45").Value = -3542.8333

$ws.Range("H63").Value = 7000
$ws.Range("I63").Value = 7000
$ws.Range("K63").Value = 7000
$ws.Range("M63").Value = -6314

$ws.Range("H66").Value = 7000
$ws.Range("I66").Value = 7000
$ws.Range("K66").Value = 35000
$ws.Range("M66").Value = -31568

$ws.Range("H102").Value = 2363.6365
$ws.Range("I102").Value = 1600
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 1600
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = 22
$ws.Range("N102").Value = -13244

$ws.Range("H122").Value = 43617.89
$ws.Range("I122").Value = 77296.2
$ws.Range("K122").Value = 231888.6
$ws.Range("M122").Value = -229438.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 42757.6
$ws.Range("I11").Value = 1886.5
$ws.Range("K11").Value = 1886.5
$ws.Range("M11").Value = -1746.5

$ws.Range("H99").Value = 1159.6154
$ws.Range("I99").Value = 1213.5
$ws.Range("K99").Value = 1213.5
$ws.Range("M99").Value = 284.5

$ws.Range("H108").Value = 90682.5
$ws.Range("J108").Value = 90682.5
$ws.Range("L108").Value = 90682.5
$ws.Range("N108").Value = -98362.5

$ws.Range("H134").Value = 10176.23
$ws.Range("I134").Value = 10690.917
$ws.Range("K134").Value = 32072.751
$ws.Range("M134").Value = -29537.751

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1047
$ws.Range("I16").Value = 1044
$ws.Range("J16").Value = 1054.5
$ws.Range("K16").Value = 1044
$ws.Range("L16").Value = 1054.5
$ws.Range("M16").Value = -757
$ws.Range("N16").Value = -1628.5

$ws.Range("H22").Value = 15626333
$ws.Range("I22").Value = 1333
$ws.Range("J22").Value = 20834666
$ws.Range("K22").Value = 1333
$ws.Range("L22").Value = 20834666
$ws.Range("M22").Value = -983
$ws.Range("N22").Value = -20835366

$ws.Range("H58").Value = 3107153.8
$ws.Range("I58").Value = 4832534.5
$ws.Range("J58").Value = 1468
$ws.Range("K58").Value = 4832534.5
$ws.Range("L58").Value = 1468
$ws.Range("M58").Value = -4832331.5
$ws.Range("N58").Value = -1874

$ws.Range("H113").Value = 1047
$ws.Range("I113").Value = 1044
$ws.Range("J113").Value = 1054.5
$ws.Range("K113").Value = 1044
$ws.Range("L113").Value = 1054.5
$ws.Range("M113").Value = 1126
$ws.Range("N113").Value = -5394.5

$ws.Range("H122").Value = 1027.5
$ws.Range("I122").Value = 1027.5
$ws.Range("K122").Value = 3082.5
$ws.Range("M122").Value = -632.5

$ws.Range("H136").Value = 3107153.8
$ws.Range("I136").Value = 4832534.5
$ws.Range("J136").Value = 1468
$ws.Range("K136").Value = 14497603.5
$ws.Range("L136").Value = 4404
$ws.Range("M136").Value = -14495053.5
$ws.Range("N136").Value = -9504

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 2073.5
$ws.Range("I46").Value = 1150
$ws.Range("J46").Value = 2997
$ws.Range("K46").Value = 3450
$ws.Range("L46").Value = 8991
$ws.Range("M46").Value = -3359
$ws.Range("N46").Value = -9173

$ws.Range("H68").Value = 1624.8298
$ws.Range("J68").Value = 1798.2565
$ws.Range("L68").Value = 5394.7695
$ws.Range("N68").Value = -7016.7695

$ws.Range("H71").Value = 1624.8298
$ws.Range("J71").Value = 1798.2565
$ws.Range("L71").Value = 16184.3085
$ws.Range("N71").Value = -24296.3085

$ws.Range("H87").Value = 10257.4
$ws.Range("I87").Value = 429
$ws.Range("K87").Value = 1287
$ws.Range("M87").Value = -39

$ws.Range("H90").Value = 10257.4
$ws.Range("I90").Value = 429
$ws.Range("K90").Value = 3861
$ws.Range("M90").Value = 2379

$ws.Range("H104").Value = 5370.4
$ws.Range("J104").Value = 5561.5557
$ws.Range("L104").Value = 16684.6671
$ws.Range("N104").Value = -21926.6671

$ws.Range("H131").Value = 8347831
$ws.Range("J131").Value = 15779.491
$ws.Range("L131").Value = 47338.473
$ws.Range("N131").Value = -57418.473

$ws.Range("H132").Value = 1038.8334
$ws.Range("I132").Value = 910
$ws.Range("J132").Value = 1054.9375
$ws.Range("K132").Value = 8190
$ws.Range("L132").Value = 9494.4375
$ws.Range("M132").Value = -5660
$ws.Range("N132").Value = -14554.4375

$ws.Range("H137").Value = 3944.1667
$ws.Range("I137").Value = 2024
$ws.Range("J137").Value = 4682.6924
$ws.Range("K137").Value = 6072
$ws.Range("L137").Value = 14048.0772
$ws.Range("M137").Value = -972
$ws.Range("N137").Value = -24248.0772

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2423.9473
$ws.Range("I122").Value = 2051.25
$ws.Range("K122").Value = 6153.75
$ws.Range("M122").Value = -3703.75

$ws.Range("H123").Value = 10326
$ws.Range("J123").Value = 10326
$ws.Range("L123").Value = 10326
$ws.Range("N123").Value = -15226

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3552.6667
$ws.Range("I7").Value = 2557.2856
$ws.Range("K7").Value = 2557.2856
$ws.Range("M7").Value = -2445.2856

$ws.Range("H122").Value = 10259.2
$ws.Range("I122").Value = 8740
$ws.Range("K122").Value = 26220
$ws.Range("M122").Value = -23770

$ws.Range("H126").Value = 3552.6667
$ws.Range("I126").Value = 2557.2856
$ws.Range("K126").Value = 7671.8568
$ws.Range("M126").Value = -5201.8568

$ws.Range("H132").Value = 3745.5386
$ws.Range("I132").Value = 1039.6
$ws.Range("K132").Value = 3118.8
$ws.Range("M132").Value = -588.7999999999997

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1774.1
$ws.Range("I81").Value = 1749
$ws.Range("K81").Value = 3498
$ws.Range("M81").Value = -2437

$ws.Range("H84").Value = 1774.1
$ws.Range("I84").Value = 1749
$ws.Range("K84").Value = 17490
$ws.Range("M84").Value = -12186

$ws.Range("H132").Value = 1764.6945
$ws.Range("I132").Value = 1426.3334
$ws.Range("J132").Value = 2441.4167
$ws.Range("K132").Value = 4279.0002
$ws.Range("L132").Value = 7324.250100000001
$ws.Range("M132").Value = -1749.0002
$ws.Range("N132").Value = -12384.2501
